$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44335
$ws.Range("J2").Value = 150
$ws.Range("M2").Value = 633
$ws.Range("P2").Value = 106

# Row 3
$ws.Range("D3").Value2 = 44335

# Row 4
$ws.Range("D4").Value2 = 44230

# Row 5
$ws.Range("D5").Value2 = 44230

# Row 6
$ws.Range("D6").Value2 = 44293
$ws.Range("J6").Value = 100

# Row 7
$ws.Range("D7").Value2 = 44293
$ws.Range("J7").Value = 50

# Row 8
$ws.Range("D8").Value2 = 44358
$ws.Range("J8").Value = 200

# Row 9
$ws.Range("D9").Value2 = 44358
$ws.Range("J9").Value = 100

# Row 10
$ws.Range("D10").Value2 = 44525

# Row 11
$ws.Range("D11").Value2 = 44525

# Row 12
$ws.Range("D12").Value2 = 44328
$ws.Range("J12").Value = 100

# Row 13
$ws.Range("D13").Value2 = 44328
$ws.Range("J13").Value = 50

# Row 14
$ws.Range("D14").Value2 = 44491
$ws.Range("J14").Value = 200
$ws.Range("O14").Value = "Región Metropolitana"

# Row 15
$ws.Range("D15").Value2 = 44491
$ws.Range("J15").Value = 100
$ws.Range("O15").Value = "Región Metropolitana"

# Row 16
$ws.Range("D16").Value2 = 44321
$ws.Range("J16").Value = 100

# Row 17
$ws.Range("D17").Value2 = 44321
$ws.Range("J17").Value = 50

# Row 18
$ws.Range("D18").Value2 = 44308
$ws.Range("J18").Value = 200
$ws.Range("M18").Value = 650
$ws.Range("P18").Value = 108

# Row 19
$ws.Range("D19").Value2 = 44308
$ws.Range("J19").Value = 100

# Row 20
$ws.Range("D20").Value2 = 44188
$ws.Range("O20").Value = "Región de Ñuble"

# Row 21
$ws.Range("D21").Value2 = 44188
$ws.Range("O21").Value = "Región de Ñuble"
